# Update КР-2 ("F") scores and the corresponding "Примечания к КР-2" (N:S,
# merged) notes for two students whose work was re-checked:
#   - row 22: Титова Надежда Алексеевна      -> КР-2 score 3 -> 4
#   - row 23: Хрищанович Полина Чеславовна   -> КР-2 score 2 -> 4
# Both rows get the note "переписаны верно все номера" (already used
# verbatim elsewhere in the sheet, e.g. row 10) in their merged N:S cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F22").Value = 4
$ws.Range("N22").Value = "переписаны верно все номера"

$ws.Range("F23").Value = 4
$ws.Range("N23").Value = "переписаны верно все номера"

# Leave the same cell selected as in the authored workbook.
$ws.Range("N24:S24").Select()
